$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update SalesQuantity (K) and Turnover (L) for the Nutella row (row 4)
$ws.Range("K4").Value = 51
$ws.Range("L4").Value = 133.85

# Update the totals row (row 5)
$ws.Range("K5").Value = 70
$ws.Range("L5").Value = 172.68
